# Commit: "Done Creat data and Update Test"
#
# 1. Rename the shopping-cart related sheets to their new, clearer names.
# 2. Set a (blank-ish) value on CartAdd!A3 that previously had no value.
# 3. Move the "active/selected" tab from ChangePass to CartDelete
#    (this updates both the workbook-level activeTab and the
#    sheet-level tabSelected flag).

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------
$wb.Worksheets.Item("SCart").Name    = "Cart"
$wb.Worksheets.Item("SCartAdd").Name = "CartAdd"
$wb.Worksheets.Item("SCartUpd").Name = "CartUpdate"
$wb.Worksheets.Item("SCartDel").Name = "CartDelete"

# --- 2. New cell value on CartAdd!A3 -----------------------------------
$wsCartAdd = $wb.Worksheets.Item("CartAdd")
$wsCartAdd.Range("A3").Value = "   "

# --- 3. Move the active tab to CartDelete ------------------------------
$wsCartDelete = $wb.Worksheets.Item("CartDelete")
$wsCartDelete.Activate()
